# Trading update: 2026-02-17 08:59:09
# Appends a new (still OPEN) MarketMaking trade as row 83 to both the
# "All Trades" and "MarketMaking" sheets, which previously ended at row 82
# (A1:Q82 -> A1:Q83).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")
$row = 83

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item($row, 1).Value = 82                # Trade #

    # Force text storage so the ISO-looking date string isn't silently
    # auto-parsed into a date serial number; reset the style back to the
    # sheet's default afterwards so no stray number-format style lingers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"       # Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:59:08"         # Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"             # Side
    $ws.Cells.Item($row, 6).Value = 0.95               # Entry Price

    # Exit Price: trade is still OPEN, so this stays blank.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"             # Status
    $ws.Cells.Item($row, 9).Value = 0                  # P&L %
    $ws.Cells.Item($row, 10).Value = 0                 # P&L $
    $ws.Cells.Item($row, 11).Value = 100.503534388353  # Capital After
    $ws.Cells.Item($row, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason

    # Exit Reason: trade is still OPEN, so this stays blank.
    $ws.Cells.Item($row, 16).NumberFormat = "@"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 16).Style = "Normal"

    $ws.Cells.Item($row, 17).Value = 0                 # Duration (min)
}
